$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title shape ("A" + " " + "slide"): consolidate the separate text runs
# into a single run by clearing the range and re-inserting the full text,
# which lets PowerPoint collapse it to one <a:r> instead of three.
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Delete() | Out-Null
$titleRange.InsertAfter("A slide") | Out-Null

# Caption textbox ("Followed" + " " + "by" + " " + "a" + " " + "picture"):
# same consolidation down to a single run.
$captionRange = $s.Shapes.Item(4).TextFrame.TextRange
$captionRange.Delete() | Out-Null
$captionRange.InsertAfter("Followed by a picture") | Out-Null
